# 10 samples from Hawaii 10.6.2019
# Append two new rows (56, 57) of CRM accuracy data below the existing
# table, following the same pattern used by the previous rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 55 into the two new rows so that the
# date style (s="1") and general row look carries over without creating
# any new cell styles.
$ws.Range("A55:F55").Copy($ws.Range("A56:F56"))
$ws.Range("A55:F55").Copy($ws.Range("A57:F57"))

# --- Row 56: Raw TA sample ---
$ws.Range("A56").Value = 43744
$ws.Range("B56").Value = 2217.7131926780899
$ws.Range("C56").Value = 2207.0300000000002
$ws.Range("D56").Formula = "=100*(B56-C56)/C56"
$ws.Range("E56").Value = 169
$ws.Range("F56").Value = "Raw TA;opened crm (10/5/2019"

# --- Row 57: TA Evap sample ---
$ws.Range("A57").Value = 43744
$ws.Range("B57").Value = 2210.1355849582301
$ws.Range("C57").Value = 2207.0300000000002
$ws.Range("D57").Formula = "=100*(B57-C57)/C57"
$ws.Range("E57").Value = 169
$ws.Range("F57").Value = "TA Evap;opened crm (10/5/2019"

# Scroll/selection state, matching where the author ended up after adding
# the new rows (cosmetic view state only).
$ws.Range("E58").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
